$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 12
$ws.Range("H12").Value = 190.16667
$ws.Range("I12").Value = 190.16667
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 190.16667
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -20.16667000000001
$ws.Range("N12").ClearContents()
# Row 98
$ws.Range("H98").Value = 362727.78
$ws.Range("I98").Value = 560382.1
$ws.Range("J98").Value = 3356.182
$ws.Range("K98").Value = 560382.1
$ws.Range("L98").Value = 3356.182
$ws.Range("M98").Value = -558884.1
$ws.Range("N98").Value = -6352.182
# Row 122
$ws.Range("H122").Value = 362727.78
$ws.Range("I122").Value = 560382.1
$ws.Range("J122").Value = 3356.182
$ws.Range("K122").Value = 1681146.3
$ws.Range("L122").Value = 10068.546
$ws.Range("M122").Value = -1678696.3
$ws.Range("N122").Value = -14968.546
# Row 132
$ws.Range("H132").Value = 361078.9
$ws.Range("I132").Value = 434884.16
$ws.Range("K132").Value = 1304652.48
$ws.Range("M132").Value = -1302122.48
# Row 135
$ws.Range("H135").Value = 6589.3
$ws.Range("I135").Value = 7229.1113
$ws.Range("K135").Value = 65062.00169999999
$ws.Range("M135").Value = -62527.00169999999

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 14094.506
$ws.Range("I32").Value = 1431.0448
$ws.Range("J32").Value = 56517.1
$ws.Range("K32").Value = 1431.0448
$ws.Range("L32").Value = 56517.1
$ws.Range("M32").Value = -1144.0448
$ws.Range("N32").Value = -57091.1
# Row 74
$ws.Range("H74").Value = 10574
$ws.Range("I74").Value = 2156.6155
$ws.Range("K74").Value = 2156.6155
$ws.Range("M74").Value = -1282.6155
# Row 77
$ws.Range("H77").Value = 10574
$ws.Range("I77").Value = 2156.6155
$ws.Range("K77").Value = 10783.0775
$ws.Range("M77").Value = -6415.077499999999
# Row 122
$ws.Range("H122").Value = 1465.8
$ws.Range("I122").Value = 1109.6666
$ws.Range("K122").Value = 3328.9998
$ws.Range("M122").Value = -878.9998000000001
# Row 133
$ws.Range("H133").Value = 32500
$ws.Range("J133").Value = 32500
$ws.Range("L133").Value = 32500
$ws.Range("N133").Value = -37560
# Row 139
$ws.Range("H139").Value = 37750
$ws.Range("J139").Value = 37750
$ws.Range("L139").Value = 37750
$ws.Range("N139").Value = -48030

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 1786.55
$ws.Range("I20").Value = 1721.1818
$ws.Range("J20").Value = 1866.4445
$ws.Range("K20").Value = 1721.1818
$ws.Range("L20").Value = 1866.4445
$ws.Range("M20").Value = -1474.1818
$ws.Range("N20").Value = -2360.4445

$ws = $wb.Worksheets.Item("CRP")
# Row 13
$ws.Range("H13").Value = 172566.67
$ws.Range("J13").Value = 258750
$ws.Range("L13").Value = 258750
$ws.Range("N13").Value = -259028
# Row 58
$ws.Range("H58").Value = 2607.762
$ws.Range("I58").Value = 1833.1666
$ws.Range("J58").Value = 3640.5557
$ws.Range("K58").Value = 1833.1666
$ws.Range("L58").Value = 3640.5557
$ws.Range("M58").Value = -1630.1666
$ws.Range("N58").Value = -4046.5557
# Row 134
$ws.Range("H134").Value = 2900.6
$ws.Range("I134").Value = 1430.0526
$ws.Range("J134").Value = 7557.3335
$ws.Range("K134").Value = 4290.1578
$ws.Range("L134").Value = 22672.0005
$ws.Range("M134").Value = -1755.1578
$ws.Range("N134").Value = -27742.0005
# Row 136
$ws.Range("H136").Value = 2607.762
$ws.Range("I136").Value = 1833.1666
$ws.Range("J136").Value = 3640.5557
$ws.Range("K136").Value = 5499.4998
$ws.Range("L136").Value = 10921.6671
$ws.Range("M136").Value = -2949.4998
$ws.Range("N136").Value = -16021.6671

$ws = $wb.Worksheets.Item("CUL")
# Row 33
$ws.Range("H33").Value = 315.13333
$ws.Range("I33").Value = 271.14285
$ws.Range("J33").Value = 353.625
$ws.Range("K33").Value = 1626.8571
$ws.Range("L33").Value = 2121.75
$ws.Range("M33").Value = -1343.8571
$ws.Range("N33").Value = -2687.75
# Row 131
$ws.Range("H131").Value = 1653.6818
$ws.Range("J131").Value = 1994.5
$ws.Range("L131").Value = 5983.5
$ws.Range("N131").Value = -16063.5

$ws = $wb.Worksheets.Item("GSM")
# Row 7
$ws.Range("H7").Value = 3000
$ws.Range("J7").Value = 3000
$ws.Range("L7").Value = 3000
$ws.Range("N7").Value = -3224
# Row 8
$ws.Range("H8").Value = 3000
$ws.Range("J8").Value = 3000
$ws.Range("L8").Value = 3000
$ws.Range("N8").Value = -3278
# Row 122
$ws.Range("H122").Value = 1390251.2
$ws.Range("I122").Value = 1853001.6
$ws.Range("K122").Value = 5559004.800000001
$ws.Range("M122").Value = -5556554.800000001
# Row 139
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 3
$ws.Range("H3").Value = 3000
$ws.Range("J3").Value = 3000
$ws.Range("L3").Value = 3000
$ws.Range("N3").Value = -3224
# Row 15
$ws.Range("H15").Value = 3000
$ws.Range("J15").Value = 3000
$ws.Range("L15").Value = 3000
$ws.Range("N15").Value = -3340
# Row 20
$ws.Range("H20").Value = 5000
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()
# Row 21
$ws.Range("H21").Value = 5000
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 5000
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 5000
$ws.Range("M21").ClearContents()
$ws.Range("N21").Value = -5348
# Row 24
$ws.Range("H24").Value = 26500
$ws.Range("I24").Value = 50000
$ws.Range("K24").Value = 50000
$ws.Range("M24").Value = -49657
# Row 136
$ws.Range("H136").Value = 9020
$ws.Range("I136").Value = 5319.6
$ws.Range("J136").Value = 11663.143
$ws.Range("K136").Value = 15958.8
$ws.Range("L136").Value = 34989.429
$ws.Range("M136").Value = -13408.8
$ws.Range("N136").Value = -40089.429

$ws = $wb.Worksheets.Item("WVR")
# Row 17
$ws.Range("H17").Value = 2749.25
$ws.Range("J17").Value = 3000
$ws.Range("L17").Value = 3000
$ws.Range("N17").Value = -3344
# Row 20
$ws.Range("H20").Value = 5499.75
$ws.Range("I20").Value = 999.5
$ws.Range("J20").Value = 10000
$ws.Range("K20").Value = 999.5
$ws.Range("L20").Value = 10000
$ws.Range("M20").Value = -759.5
$ws.Range("N20").Value = -10480
# Row 122
$ws.Range("H122").Value = 79040
$ws.Range("I122").Value = 201120.8
$ws.Range("J122").Value = 2739.5
$ws.Range("K122").Value = 603362.3999999999
$ws.Range("L122").Value = 8218.5
$ws.Range("M122").Value = -600912.3999999999
$ws.Range("N122").Value = -13118.5
